$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.67803479114983412
$ws.Range("C2").Value = 0.54150621371675611
$ws.Range("D2").Value = 0.97288749503397409
$ws.Range("F4").Value = 0.95273202936010803
$ws.Range("C5").Value = 0.88200444055247451
$ws.Range("L5").Value = 0.90557777691126151
$ws.Range("E6").Value = 0.97124550683980493
$ws.Range("G6").Value = 0.83878787648448683
$ws.Range("H6").Value = 0.96577833579014993
$ws.Range("AH6").Value = 0.86377972732039421
$ws.Range("E7").Value = 0.87558583060844009
$ws.Range("C8").Value = 0.92037825725492484
$ws.Range("G8").Value = 0.92200598780406762
$ws.Range("AU8").Value = 0.93777673633464498
$ws.Range("G9").Value = 0.90056726492772987
$ws.Range("J9").Value = 0.64214449778127114
$ws.Range("H10").Value = 0.79162633636021129
$ws.Range("L10").Value = 0.88600758086320652
$ws.Range("AX10").Value = 0.80374068561669831
$ws.Range("BH10").Value = 0.82265167708283382
$ws.Range("I11").Value = 0.99048863995723546
$ws.Range("J11").Value = 0.7107858083945835
$ws.Range("K12").Value = 0.80642629184520076
$ws.Range("M12").Value = 0.55637028250836917
$ws.Range("K13").Value = 0.71563985981266964
$ws.Range("O13").Value = 0.61918823160932579
$ws.Range("M14").Value = 0.94019250393537179
$ws.Range("N15").Value = 0.88093750698839091
$ws.Range("Q15").Value = 0.95884603683489489
$ws.Range("AY15").Value = 0.96000543082232137
$ws.Range("O16").Value = 0.85737413787442618
$ws.Range("Q16").Value = 0.5350429489478209
$ws.Range("R16").Value = 0.98355729400800107
$ws.Range("R17").Value = 0.99936229994592696
$ws.Range("H19").Value = 0.85852537995773359
$ws.Range("Q19").Value = 0.75047496476156972
$ws.Range("R19").Value = 0.83355280616886152
$ws.Range("T19").Value = 0.9190567380578899
$ws.Range("R20").Value = 0.80370193966072456
$ws.Range("N21").Value = 0.76360999904246118
$ws.Range("S21").Value = 0.95255422392249145
$ws.Range("T21").Value = 0.83319775947561592
$ws.Range("V21").Value = 0.83083712574213864
$ws.Range("T22").Value = 0.93604060545225798
$ws.Range("W22").Value = 0.95104229525100825
$ws.Range("V24").Value = 0.83729147107772561
$ws.Range("Z24").Value = 0.7764972605777547
$ws.Range("AZ24").Value = 0.81078121413340642
$ws.Range("W25").Value = 0.72263166818914237
$ws.Range("Y26").Value = 0.72510601481172987
$ws.Range("AA26").Value = 0.89505782868874872
$ws.Range("AB26").Value = 0.78782655330303175
$ws.Range("Y27").Value = 0.97643346860326075
$ws.Range("AA28").Value = 0.81971758955882934
$ws.Range("AC28").Value = 0.88844864534379742
$ws.Range("X29").Value = 0.7291483854405314
$ws.Range("AA29").Value = 0.68434956475421882
$ws.Range("AD29").Value = 0.85539019028888053
$ws.Range("AB30").Value = 0.82451925973670392
$ws.Range("AE30").Value = 0.81884797596899106
$ws.Range("N31").Value = 0.98681131968814095
$ws.Range("AF31").Value = 0.67790647497535472
$ws.Range("AG31").Value = 0.95749742473639343
$ws.Range("AF33").Value = 0.94599866184163162
$ws.Range("BM33").Value = 0.99486495228444261
$ws.Range("AF34").Value = 0.69481494668902166
$ws.Range("AG34").Value = 0.69708130146717662
$ws.Range("AS34").Value = 0.86107138554654628
$ws.Range("AH35").Value = 0.72356885342728994
$ws.Range("AK36").Value = 0.90249539850787741
$ws.Range("BI36").Value = 0.90448437358814782
$ws.Range("BM36").Value = 0.86304587135711119
$ws.Range("AM37").Value = 0.91899484689294775
$ws.Range("BE37").Value = 0.94588533245969686
$ws.Range("AK38").Value = 0.848761380938995
$ws.Range("AM38").Value = 0.89930443638345314
$ws.Range("AN38").Value = 0.52127756934793257
$ws.Range("AI39").Value = 0.80209955144346634
$ws.Range("W40").Value = 0.82471752458801673
$ws.Range("AM40").Value = 0.93368560821275048
$ws.Range("AM41").Value = 0.95154558136803535
$ws.Range("AN41").Value = 0.86252758825440479
$ws.Range("BH41").Value = 0.65695231980018742
$ws.Range("AN42").Value = 0.94794085771131464
$ws.Range("AO42").Value = 0.93100034202679494
$ws.Range("AQ42").Value = 0.66864853068312802
$ws.Range("AR43").Value = 0.95099495242011711
$ws.Range("BH43").Value = 0.82652350780444428
$ws.Range("AP44").Value = 0.91474042081035889
$ws.Range("AQ45").Value = 0.9829772245583217
$ws.Range("AR45").Value = 0.96890728593019415
$ws.Range("AR46").Value = 0.70162698430614312
$ws.Range("AS46").Value = 0.76913474449232822
$ws.Range("AT47").Value = 0.99334594746585747
$ws.Range("AW47").Value = 0.67281940069482582
$ws.Range("AT48").Value = 0.58215015269244086
$ws.Range("AV49").Value = 0.67353648712874992
$ws.Range("AX49").Value = 0.97326508207775597
$ws.Range("BM49").Value = 0.68209840043648207
$ws.Range("AV50").Value = 0.98891242609354979
$ws.Range("AZ50").Value = 0.96151681693828794
$ws.Range("BG50").Value = 0.89572529493956377
$ws.Range("BA51").Value = 0.83483212012531616
$ws.Range("AI52").Value = 0.5875295190472124
$ws.Range("BA52").Value = 0.85067070819912938
$ws.Range("BB52").Value = 0.98290981255491949
$ws.Range("BC53").Value = 0.901961027934989
$ws.Range("BA54").Value = 0.67941281477392179
$ws.Range("BE54").Value = 0.99084464045993581
$ws.Range("BB55").Value = 0.70262328142551933
$ws.Range("BD55").Value = 0.96358242357875812
$ws.Range("BB56").Value = 0.71977807724890963
$ws.Range("BF56").Value = 0.83834761909553401
$ws.Range("BC57").Value = 0.71652098056573177
$ws.Range("BD57").Value = 0.90277788659312885
$ws.Range("BH58").Value = 0.91503808136051434
$ws.Range("BF59").Value = 0.96251542966482506
$ws.Range("BJ60").Value = 0.96531579577547921
$ws.Range("BG61").Value = 0.97536184198809617
$ws.Range("BJ61").Value = 0.86054725989946579
$ws.Range("J62").Value = 0.95083889959931911
$ws.Range("B63").Value = 0.83498329055797293
$ws.Range("BM63").Value = 0.84817692346122486
$ws.Range("A64").Value = 0.85169625792419446
$ws.Range("AD64").Value = 0.58674261559540275
$ws.Range("AA65").Value = 0.9318153320876057
$ws.Range("BL65").Value = 0.69401642520222717
$ws.Range("D66").Value = 0.97963536766872794
$ws.Range("J66").Value = 0.80858962243266341
$ws.Range("BL66").Value = 0.89215293854546607
$ws.Range("BM66").Value = 0.83521667918063969
$ws.Range("BP66").Value = 0.74171997896532815
$ws.Range("J67").Value = 0.82138820442584071
$ws.Range("BP67").Value = 0.96907482041076309
$ws.Range("A68").Value = 0.63238729277621542
$ws.Range("B68").Value = 0.70933300472642835
